# aggiornamento fino a 6/03
# Append three new daily rows (245-247) to the data table, continuing the
# date sequence in column A and zero values in columns B, C, D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data ends at row 244 (date serial 44318 = 2021-05-02).
# New rows continue the sequence through 2021-05-05.
$dates = @(44319, 44320, 44321)
$startRow = 245

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $prevRow = $row - 1

    # Carry forward the formatting of the row above (date style incl.
    # border/alignment/number format) instead of re-creating it, so the
    # style table stays identical to the source row.
    $ws.Range("A$prevRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

$excel.CutCopyMode = 0
